$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-01 03:16:23"

$wsZhCn.Range("H2").Value = "2016-09-01 03:16:19"
$wsZhCn.Range("K2").Value = "2016-09-01 03:16:49"

$wsDeDe.Range("H2").Value = "2016-09-01 03:16:23"
$wsDeDe.Range("K2").Value = "2016-09-01 03:16:56"
